# Scheduled-runner market data refresh: update currentAveragePrice /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ figures (cols H-N)
# for a handful of rows across all 8 crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 982.8
$ws.Range("I96").Value = 1084
$ws.Range("J96").Value = 746.6667
$ws.Range("K96").Value = 3252
$ws.Range("L96").Value = 2240.0001
$ws.Range("M96").Value = -1879
$ws.Range("N96").Value = -4986.0001

$ws.Range("H101").Value = 828
$ws.Range("I101").Value = 895.625
$ws.Range("J101").Value = 557.5
$ws.Range("K101").Value = 2686.875
$ws.Range("L101").Value = 1672.5
$ws.Range("M101").Value = -1064.875
$ws.Range("N101").Value = -4916.5

$ws.Range("H121").Value = 1500
$ws.Range("J121").Value = 1500
$ws.Range("L121").Value = 4500
$ws.Range("N121").Value = -7994

$ws.Range("H132").Value = 294912.2
$ws.Range("I132").Value = 743.80646
$ws.Range("K132").Value = 2231.41938
$ws.Range("M132").Value = 298.5806199999997

$ws.Range("H135").Value = 7351.467
$ws.Range("I135").Value = 713.3077
$ws.Range("J135").Value = 50499.5
$ws.Range("K135").Value = 6419.7693
$ws.Range("L135").Value = 454495.5
$ws.Range("M135").Value = -3884.7693
$ws.Range("N135").Value = -459565.5

$ws.Range("H137").Value = 4010.9167
$ws.Range("I137").Value = 2171
$ws.Range("J137").Value = 7690.75
$ws.Range("K137").Value = 6513
$ws.Range("L137").Value = 23072.25
$ws.Range("M137").Value = -3963
$ws.Range("N137").Value = -28172.25

$ws.Range("H138").Value = 3873.17
$ws.Range("I138").Value = 2067.682
$ws.Range("J138").Value = 5154.484
$ws.Range("K138").Value = 6203.045999999999
$ws.Range("L138").Value = 15463.452
$ws.Range("M138").Value = -1063.045999999999
$ws.Range("N138").Value = -25743.452


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4184.974
$ws.Range("I32").Value = 1962.9054
$ws.Range("K32").Value = 1962.9054
$ws.Range("M32").Value = -1675.9054

$ws.Range("H132").Value = 1487.8096
$ws.Range("I132").Value = 1481
$ws.Range("J132").Value = 1576.3334
$ws.Range("K132").Value = 4443
$ws.Range("L132").Value = 4729.0002
$ws.Range("M132").Value = -1913
$ws.Range("N132").Value = -9789.0002


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 14997.333
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("M54").Value = -4516

$ws.Range("H102").Value = 1249.8
$ws.Range("I102").Value = 1249.8
$ws.Range("K102").Value = 1249.8
$ws.Range("M102").Value = 1995.2

$ws.Range("H134").Value = 2388.1892
$ws.Range("I134").Value = 2276.0857
$ws.Range("J134").Value = 4350
$ws.Range("K134").Value = 6828.257100000001
$ws.Range("L134").Value = 13050
$ws.Range("M134").Value = -4293.257100000001
$ws.Range("N134").Value = -18120


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5914.349
$ws.Range("I31").Value = 3823.9697
$ws.Range("J31").Value = 8213.767
$ws.Range("K31").Value = 3823.9697
$ws.Range("L31").Value = 8213.767
$ws.Range("M31").Value = -3528.9697
$ws.Range("N31").Value = -8803.767

$ws.Range("H34").Value = 5914.349
$ws.Range("I34").Value = 3823.9697
$ws.Range("J34").Value = 8213.767
$ws.Range("K34").Value = 3823.9697
$ws.Range("L34").Value = 8213.767
$ws.Range("M34").Value = -3621.9697
$ws.Range("N34").Value = -8617.767

$ws.Range("H58").Value = 2824.6316
$ws.Range("I58").Value = 2881.6667
$ws.Range("K58").Value = 2881.6667
$ws.Range("M58").Value = -2678.6667

$ws.Range("H99").Value = 2805.6
$ws.Range("I99").Value = 2341.8333
$ws.Range("K99").Value = 2341.8333
$ws.Range("M99").Value = -843.8332999999998

$ws.Range("H103").Value = 22704
$ws.Range("I103").Value = 12380
$ws.Range("K103").Value = 12380
$ws.Range("M103").Value = -11208

$ws.Range("H126").Value = 2805.6
$ws.Range("I126").Value = 2341.8333
$ws.Range("K126").Value = 7025.499899999999
$ws.Range("M126").Value = -4555.499899999999

$ws.Range("H132").Value = 1178.0488
$ws.Range("I132").Value = 1097.4138
$ws.Range("J132").Value = 1372.9166
$ws.Range("K132").Value = 3292.2414
$ws.Range("L132").Value = 4118.7498
$ws.Range("M132").Value = -762.2413999999999
$ws.Range("N132").Value = -9178.7498

$ws.Range("H136").Value = 2824.6316
$ws.Range("I136").Value = 2881.6667
$ws.Range("K136").Value = 8645.000100000001
$ws.Range("M136").Value = -6095.000100000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2925.0908
$ws.Range("I137").Value = 2823
$ws.Range("J137").Value = 3047.6
$ws.Range("K137").Value = 8469
$ws.Range("L137").Value = 9142.799999999999
$ws.Range("M137").Value = -3369
$ws.Range("N137").Value = -19342.8


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 39999
$ws.Range("J130").Value = 39999
$ws.Range("L130").Value = 39999
$ws.Range("N130").Value = -50039

$ws.Range("H132").Value = 1956.55
$ws.Range("I132").Value = 1956.55
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5869.65
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -3339.65


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1796.2
$ws.Range("I22").Value = 1185.5714
$ws.Range("J22").Value = 3221
$ws.Range("K22").Value = 1185.5714
$ws.Range("L22").Value = 3221
$ws.Range("M22").Value = -890.5714
$ws.Range("N22").Value = -3811

$ws.Range("H27").Value = 1796.2
$ws.Range("I27").Value = 1185.5714
$ws.Range("J27").Value = 3221
$ws.Range("K27").Value = 1185.5714
$ws.Range("L27").Value = 3221
$ws.Range("M27").Value = -1078.5714
$ws.Range("N27").Value = -3435

$ws.Range("H93").Value = 1650.0625
$ws.Range("I93").Value = 2651.75
$ws.Range("J93").Value = 648.375
$ws.Range("K93").Value = 2651.75
$ws.Range("L93").Value = 648.375
$ws.Range("M93").Value = -1403.75
$ws.Range("N93").Value = -3144.375

$ws.Range("H104").Value = 13898
$ws.Range("J104").Value = 13898
$ws.Range("L104").Value = 13898
$ws.Range("N104").Value = -20886

$ws.Range("H130").Value = 55999.5
$ws.Range("J130").Value = 55999.5
$ws.Range("L130").Value = 55999.5
$ws.Range("N130").Value = -66039.5

$ws.Range("H132").Value = 1953.6578
$ws.Range("I132").Value = 1737.6897
$ws.Range("J132").Value = 2649.5557
$ws.Range("K132").Value = 5213.0691
$ws.Range("L132").Value = 7948.6671
$ws.Range("M132").Value = -2683.0691
$ws.Range("N132").Value = -13008.6671

$ws.Range("H136").Value = 1091.6428
$ws.Range("I136").Value = 928.3
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 2784.9
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -234.8999999999996
$ws.Range("N136").Value = -9600


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3032981.2
$ws.Range("I96").Value = 4042973.8
$ws.Range("J96").Value = 3004
$ws.Range("K96").Value = 4042973.8
$ws.Range("L96").Value = 3004
$ws.Range("M96").Value = -4041600.8
$ws.Range("N96").Value = -5750

$ws.Range("H132").Value = 2299.0667
$ws.Range("I132").Value = 2283.182
$ws.Range("J132").Value = 2342.75
$ws.Range("K132").Value = 6849.545999999999
$ws.Range("L132").Value = 7028.25
$ws.Range("M132").Value = -4319.545999999999
$ws.Range("N132").Value = -12088.25


# Remove cell that no longer has a value in the source row (GSM row 132, column M)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M132").ClearContents()
